$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6747177243232727
$ws.Range("B1").Value = 0.9657843708992004
$ws.Range("C1").Value = 4.754477500915527
$ws.Range("D1").Value = 1.929329037666321
$ws.Range("E1").Value = 1.163437366485596
